$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 27783044
$ws.Range("I70").Value = 2612.25
$ws.Range("J70").Value = 50007388
$ws.Range("K70").Value = 7836.75
$ws.Range("L70").Value = 150022164
$ws.Range("M70").Value = -7566.75
$ws.Range("N70").Value = -150022704
$ws.Range("H73").Value = 27783044
$ws.Range("I73").Value = 2612.25
$ws.Range("J73").Value = 50007388
$ws.Range("K73").Value = 7836.75
$ws.Range("L73").Value = 150022164
$ws.Range("M73").Value = -6900.75
$ws.Range("N73").Value = -150024036
$ws.Range("H98").Value = 1752.6111
$ws.Range("I98").Value = 825.2143
$ws.Range("K98").Value = 825.2143
$ws.Range("M98").Value = 672.7857
$ws.Range("H116").Value = 71438696
$ws.Range("J116").Value = 9249
$ws.Range("L116").Value = 9249
$ws.Range("N116").Value = -16133
$ws.Range("H122").Value = 1752.6111
$ws.Range("I122").Value = 825.2143
$ws.Range("K122").Value = 2475.6429
$ws.Range("M122").Value = -25.64289999999983
$ws.Range("H138").Value = 2435.28
$ws.Range("I138").Value = 916.875
$ws.Range("J138").Value = 2724.5
$ws.Range("K138").Value = 2750.625
$ws.Range("L138").Value = 8173.5
$ws.Range("M138").Value = 2389.375
$ws.Range("N138").Value = -18453.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2090
$ws.Range("I2").Value = 1913.8667
$ws.Range("J2").Value = 2255.125
$ws.Range("K2").Value = 1913.8667
$ws.Range("L2").Value = 2255.125
$ws.Range("M2").Value = -1800.8667
$ws.Range("N2").Value = -2481.125
$ws.Range("H116").Value = 2090
$ws.Range("I116").Value = 1913.8667
$ws.Range("J116").Value = 2255.125
$ws.Range("K116").Value = 1913.8667
$ws.Range("L116").Value = 2255.125
$ws.Range("M116").Value = 380.1333
$ws.Range("N116").Value = -6843.125
$ws.Range("H122").Value = 3102
$ws.Range("I122").Value = 2120.5
$ws.Range("J122").Value = 4574.25
$ws.Range("K122").Value = 6361.5
$ws.Range("L122").Value = 13722.75
$ws.Range("M122").Value = -3911.5
$ws.Range("N122").Value = -18622.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2090
$ws.Range("I3").Value = 1913.8667
$ws.Range("J3").Value = 2255.125
$ws.Range("K3").Value = 1913.8667
$ws.Range("L3").Value = 2255.125
$ws.Range("M3").Value = -1799.8667
$ws.Range("N3").Value = -2483.125
$ws.Range("H86").Value = 2111.7778
$ws.Range("I86").Value = 1983.3334
$ws.Range("J86").Value = 2368.6667
$ws.Range("K86").Value = 1983.3334
$ws.Range("L86").Value = 2368.6667
$ws.Range("M86").Value = -860.3334
$ws.Range("N86").Value = -4614.6667
$ws.Range("H89").Value = 2111.7778
$ws.Range("I89").Value = 1983.3334
$ws.Range("J89").Value = 2368.6667
$ws.Range("K89").Value = 9916.666999999999
$ws.Range("L89").Value = 11843.3335
$ws.Range("M89").Value = -4300.666999999999
$ws.Range("N89").Value = -23075.3335
$ws.Range("H99").Value = 3430.4375
$ws.Range("I99").Value = 2948
$ws.Range("J99").Value = 3719.9
$ws.Range("K99").Value = 2948
$ws.Range("L99").Value = 3719.9
$ws.Range("M99").Value = -1450
$ws.Range("N99").Value = -6715.9
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1389.6086
$ws.Range("I107").Value = 681.7778
$ws.Range("K107").Value = 681.7778
$ws.Range("M107").Value = 1238.2222
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 4127
$ws.Range("I18").Value = 3056.8572
$ws.Range("K18").Value = 9170.571599999999
$ws.Range("M18").Value = -9001.571599999999
$ws.Range("H56").Value = 8828.571
$ws.Range("I56").Value = 8828.571
$ws.Range("K56").Value = 8828.571
$ws.Range("M56").Value = -8298.571
$ws.Range("H62").Value = 5222
$ws.Range("J62").Value = 4963
$ws.Range("L62").Value = 14889
$ws.Range("N62").Value = -16261
$ws.Range("H65").Value = 5222
$ws.Range("J65").Value = 4963
$ws.Range("L65").Value = 44667
$ws.Range("N65").Value = -51531
$ws.Range("H113").Value = 979.4286
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 1042.6666
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 3127.9998
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -7467.9998
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 46177.4
$ws.Range("J52").Value = 45296.5
$ws.Range("L52").Value = 45296.5
$ws.Range("N52").Value = -45814.5
$ws.Range("H102").Value = 2522.2307
$ws.Range("I102").Value = 2506.0908
$ws.Range("J102").Value = 2611
$ws.Range("K102").Value = 2506.0908
$ws.Range("L102").Value = 2611
$ws.Range("M102").Value = -884.0907999999999
$ws.Range("N102").Value = -5855
$ws.Range("H141").Value = 47143
$ws.Range("J141").Value = 47143
$ws.Range("L141").Value = 47143
$ws.Range("N141").Value = -57503
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2177.6667
$ws.Range("I82").Value = 1950.1
$ws.Range("K82").Value = 1950.1
$ws.Range("M82").Value = -1589.1
$ws.Range("H85").Value = 2177.6667
$ws.Range("I85").Value = 1950.1
$ws.Range("K85").Value = 1950.1
$ws.Range("M85").Value = -702.0999999999999
$ws.Range("H136").Value = 10529.591
$ws.Range("I136").Value = 5657.9165
$ws.Range("J136").Value = 16375.6
$ws.Range("K136").Value = 16973.7495
$ws.Range("L136").Value = 49126.8
$ws.Range("M136").Value = -14423.7495
$ws.Range("N136").Value = -54226.8
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 81216.39999999999
$ws.Range("J16").Value = 81216.39999999999
$ws.Range("L16").Value = 81216.39999999999
$ws.Range("N16").Value = -81800.39999999999
$ws.Range("H49").Value = 16037.333
$ws.Range("I49").Value = 11556
$ws.Range("K49").Value = 11556
$ws.Range("M49").Value = -11326
$ws.Range("H81").Value = 50278.523
$ws.Range("I81").Value = 68663.39999999999
$ws.Range("J81").Value = 4316.3335
$ws.Range("K81").Value = 137326.8
$ws.Range("L81").Value = 8632.666999999999
$ws.Range("M81").Value = -136265.8
$ws.Range("N81").Value = -10754.667
$ws.Range("H84").Value = 50278.523
$ws.Range("I84").Value = 68663.39999999999
$ws.Range("J84").Value = 4316.3335
$ws.Range("K84").Value = 686634
$ws.Range("L84").Value = 43163.335
$ws.Range("M84").Value = -681330
$ws.Range("N84").Value = -53771.335
$ws.Range("H136").Value = 1635.1714
$ws.Range("I136").Value = 1086.6786
$ws.Range("J136").Value = 3829.1428
$ws.Range("K136").Value = 3260.0358
$ws.Range("L136").Value = 11487.4284
$ws.Range("M136").Value = -710.0357999999997
$ws.Range("N136").Value = -16587.4284
